$d = $word.ActiveDocument

# --- Locate the target paragraph -------------------------------------------------
# The document has (in reading order): the big features table, whose last
# row/cell contains "Implementing a PDF viewer ..."; then a run of blank
# paragraphs; then a paragraph that only holds a manual page break; then the
# "V3.0" heading. We need the *last* blank paragraph before that page break
# (the one that already carries spacing/contextualSpacing/jc formatting but
# has no runs yet) - that's where the new wishlist items get appended.

$count = $d.Paragraphs.Count

# 1) Find the anchor paragraph ("Implementing a PDF viewer ...").
$anchorIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -match "^Implementing a PDF viewer") {
        $anchorIdx = $i
        break
    }
}
if ($anchorIdx -eq -1) {
    throw "Could not locate anchor paragraph 'Implementing a PDF viewer...'"
}

# 2) Walk forward to find the paragraph that starts with a manual page
#    break (character code 12).
$pageBreakIdx = -1
$scanLimit = [Math]::Min($count, $anchorIdx + 15)
for ($j = $anchorIdx + 1; $j -le $scanLimit; $j++) {
    $t = $d.Paragraphs.Item($j).Range.Text
    if ($t.Length -gt 0 -and [int]$t[0] -eq 12) {
        $pageBreakIdx = $j
        break
    }
}
if ($pageBreakIdx -eq -1) {
    throw "Could not locate the page-break paragraph after the anchor"
}

# 3) The target is the paragraph immediately preceding the page break.
$targetIdx = $pageBreakIdx - 1
$target = $d.Paragraphs.Item($targetIdx)

# --- Apply the edit ----------------------------------------------------------------
# Fill the existing empty paragraph with its new text.
$target.Range.InsertAfter("events, technologies")

# Insert a new paragraph: "have option to filter all ticked ones".
$target.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($targetIdx + 1)
$p2.Range.InsertAfter("have option to filter all ticked ones")

# Insert three blank paragraphs.
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($targetIdx + 2)

$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Item($targetIdx + 3)

$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item($targetIdx + 4)

# Insert the final new paragraph: "entities - have icons beside names".
$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs.Item($targetIdx + 5)
$p6.Range.InsertAfter("entities - have icons beside names")
